$wb = $excel.ActiveWorkbook

# --- pdx_model sheet ---
# Row 2 and row 3 both referenced the same model_id "CRC0014LM", which now violates
# the unique constraint on model_id. Update them to distinct, valid model ids.
$pdxModel = $wb.Worksheets.Item("pdx_model")
$pdxModel.Range("A2").Value = "CRC0228PR"
$pdxModel.Range("A3").Value = "CRC0228PRaS"
$pdxModel.Range("A3").Select()

# --- model_validation sheet ---
# The model_id referenced here ("CRC0s228PR") was a typo'd value not matching any
# pdx_model row; point it at the real model id used above.
$modelValidation = $wb.Worksheets.Item("model_validation")
$modelValidation.Range("A3").Value = "CRC0228PRaS"
$modelValidation.Range("A3").Select()

# --- cell_model sheet ---
# Row 3 duplicated the "CRC0014LM" model_id already used on row 2, breaking the
# unique constraint; rename it to a distinct id.
$cellModel = $wb.Worksheets.Item("cell_model")
$cellModel.Range("A3").Value = "CRC0014LM_2"
$cellModel.Range("A3").Select()

# --- patient_sample sheet ---
# No cell values changed here, just move the remembered selection.
$patientSample = $wb.Worksheets.Item("patient_sample")
$patientSample.Range("T3").Select()

# --- patient sheet ---
# A3 previously held "SADASD" - it is no longer a valid/needed value, clear it
# (keeps its existing style s="4").
$patient = $wb.Worksheets.Item("patient")
$patient.Range("A3").Value = ""
# Row 4 / B4 held a duplicate "MALE" entry (with a one-off style) that is removed entirely.
$patient.Range("B4").Clear()
# Select the patient sheet last so it remains the active tab, matching the
# original file (only the selected cell on this tab changes, B4 -> A3).
$patient.Activate()
$patient.Range("A3").Select()
